$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "86.982.10"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").Value = "3.047.13"
$ws.Range("E3").Value = "  -4.59%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'207.50"
$ws.Range("E5").Value = "  -3.67%  "
$ws.Range("D6").Value = "'617.29"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "'0.356"
$ws.Range("E7").Value = "  -11.44%  "
$ws.Range("D8").Value = "'0.764"
$ws.Range("E8").Value = "  +10.39%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "3.041.17"
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("D11").Value = "'0.571"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "'0.176"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  -11.29%  "
$ws.Range("D14").Value = "'5.19"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "86.867.69"
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").Value = "3.607.95"
$ws.Range("E16").Value = "  -4.66%  "
$ws.Range("D17").Value = "'30.78"
$ws.Range("E17").Value = "  -7.18%  "
$ws.Range("D18").Value = "3.071.47"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("D19").Value = "'3.28"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'0.0000204"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("D21").Value = "'12.84"
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("D22").Value = "'411.01"
$ws.Range("E22").Value = "  -6.66%  "
$ws.Range("D23").Value = "'8.08"
$ws.Range("E23").Value = "  -6.55%  "
$ws.Range("D24").Value = "'4.71"
$ws.Range("E24").Value = "  -7.40%  "
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "'81.44"
$ws.Range("E26").Value = "  +7.74%  "
$ws.Range("D27").Value = "'11.04"
$ws.Range("E27").Value = "  -5.50%  "
$ws.Range("D28").Value = "3.223.39"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -12.73%  "
$ws.Range("D32").Value = "'7.89"
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("D33").Value = "'488.64"
$ws.Range("E33").Value = "  -9.27%  "
$ws.Range("D34").Value = "'3.55"
$ws.Range("E34").Value = "  -15.91%  "
$ws.Range("E35").Value = "  +8.94%  "
$ws.Range("D36").Value = "'6.52"
$ws.Range("E36").Value = "  -7.51%  "
$ws.Range("D37").Value = "'1.76"
$ws.Range("E37").Value = "  -6.11%  "
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'22.10"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'21.60"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  -5.52%  "
$ws.Range("D44").Value = "'145.95"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("E45").Value = "  -8.49%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'43.38"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.129"
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("D48").Value = "'0.0632"
$ws.Range("E48").Value = "  +6.22%  "
$ws.Range("D49").Value = "'156.07"
$ws.Range("E49").Value = "  -9.98%  "
$ws.Range("D50").Value = "'0.696"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -7.68%  "
